$p = $ppt.ActivePresentation

# "Lesson 7.4 Case Study -- Undefined Variables", slide 30 ("Next Steps"),
# Content Placeholder 2 -- the bullet "Do Guided Practice 7.3" needs to
# become "Do Guided Practices 7.2 and 7.3", split across two runs ("Do "
# and "Guided Practices 7.2 and 7.3") the way PowerPoint leaves it after an
# in-place retype of the back half of the line.
$slide = $p.Slides.Item(30)
$shape = $slide.Shapes.Item("Content Placeholder 2")
$tr = $shape.TextFrame.TextRange

$thirdPara = $tr.Paragraphs(3, 1)

$leadRun = $thirdPara.Runs(1, 1)
$leadRun.Text = "Do "
$leadRun.InsertAfter("Guided Practices 7.2 and 7.3") | Out-Null
